$wb = $excel.ActiveWorkbook

# New data for both sheets after the "fixed workflow" re-run:
# Column A (index) is always 0..14 (15 rows of data, rows 2-16).
# Column B (Cutoff) and Column C (Reaction_number) are the last 15 values
# of what used to be a 19-row series (i.e. the first 4 rows were dropped).

$colA = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14)

$nbrB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$nbrC = @(94,95,93,94,92,89,90,89,84,84,82,82,82,81,82)

$barB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$barC = @(581,579,582,576,578,578,575,574,577,576,575,575,573,575,571)

foreach ($sheetName in @("NBR", "BAR")) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "NBR") {
        $dataB = $nbrB
        $dataC = $nbrC
    } else {
        $dataB = $barB
        $dataC = $barC
    }

    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $colA[$i]
        $ws.Cells.Item($row, 2).Value = $dataB[$i]
        $ws.Cells.Item($row, 3).Value = $dataC[$i]
    }

    # Remove the now-unused trailing rows (old rows 17-20)
    $ws.Range("A17:C20").Clear()
}
